$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26, column A ("phone") was mistakenly stored as text ("71076783").
# Fix it so it is a real number, matching the rest of the column.
$ws.Cells.Item(26, 1).Value = 71076783

# Redeem points 76442781 100.0 -> append a new redemption row.
# The phone number is written as text (as the redemption importer does),
# so prefix with an apostrophe to stop Excel from re-interpreting the
# digit string as a number.
$ws.Cells.Item(27, 1).Value = "'76442781"
$ws.Cells.Item(27, 2).Value = 100
$ws.Cells.Item(27, 3).Value = "2025-08-18T18:06:34"
